$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared strings must be appended to the table in this order to match
# the target sharedStrings.xml ordering: Aurora, Glacier, White Night..., Snow...
$ws.Range("E5").Value = "Aurora, Guardian of the Radiant Skies"
$ws.Range("D5").Value = "Glacier, Incarnation of the Frozen Tempest"
$ws.Range("D4").Value = "White Night Dragon / Alexandrite Dragon / Krystal Dragon"
$ws.Range("C5").Value = "Snow, Wyvern of the Blizzard"

# Move the active selection to C5 (matches the saved cursor position in the diff)
$ws.Range("C5").Select()
